$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# New column header "habitat" (bold, like the other header cells)
$ws.Range("I1").Value = "habitat"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").IndentLevel = 0

# Column width for the new column
$ws.Columns.Item(9).ColumnWidth = 13.7

# Fill "freshwater" for every data row, centered like the other data columns
for ($r = 2; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Value = "freshwater"
    $cell.HorizontalAlignment = -4108
}

# Selection left by the editor
$ws.Range("N16").Select() | Out-Null
